# Linear-regression forecast update: the Gap_Growth_% / Ratio_Change_%
# columns (F, G) no longer carry stray 0 placeholders on the two footer
# rows ("Back to index" / "Please click to email us your opinion:") -
# those rows have no data, so F48:G49 should simply be blank, matching
# the rest of the non-data cells in those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F48:G49").ClearContents()
